# Weekly fruit/vegetable price update:
# Insert a new daily record at row 423 (pushing existing rows down by one),
# matching the "Hortaliza, Vega Modelo de Temuco - Zanahoria" sheet update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 423; this shifts rows 423:477 down to
# 424:478 and grows the used range / dimension automatically.
$ws.Rows.Item(423).Insert()

# Populate the newly inserted row 423 with the new record's data.
$ws.Cells.Item(423, 1).Value2 = 10
$ws.Cells.Item(423, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(423, 3).Value2 = "La Araucanía"
$ws.Cells.Item(423, 4).Value2 = 45077
$ws.Cells.Item(423, 5).Value2 = 9
$ws.Cells.Item(423, 6).Value2 = 100114013
$ws.Cells.Item(423, 7).Value2 = "Zanahoria"
$ws.Cells.Item(423, 8).Value2 = "Sin especificar"
$ws.Cells.Item(423, 9).Value2 = "Primera"
$ws.Cells.Item(423, 10).Value2 = 125
$ws.Cells.Item(423, 11).Value2 = 5000
$ws.Cells.Item(423, 12).Value2 = 5000
$ws.Cells.Item(423, 13).Value2 = 5000
$ws.Cells.Item(423, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(423, 15).Value2 = "Región de La Araucanía"
$ws.Cells.Item(423, 16).Value2 = 200
$ws.Cells.Item(423, 17).Value2 = 25
$ws.Cells.Item(423, 18).Value2 = "Hortaliza"

# Apply the same date-number format used by the other rows in column D.
$ws.Cells.Item(423, 4).NumberFormat = $ws.Cells.Item(424, 4).NumberFormat
